$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns that receive new values for rows 2-25 (H is untouched, stays 0)
$cols = @("C","D","E","F","G","I","J","K","L")

# New values per row (rows 2..25), one inner array per row, aligned with $cols
$data = @(
    @(0.05353901873847633,0.009189558092050021,0.1049200150058525,5.014851781193272,0.002676194052642398,3.666749870973575,0.266854974889938,4.277302916017447,0.2137801060446165),
    @(0.05261101080957786,0.008832666357061925,0.1045429673066511,5.015863496769143,0.002682328664329303,3.667720534125863,0.2670129980810927,4.138038142542996,0.213768216936117),
    @(0.05205713464987127,0.008613325950223683,0.1043491347591861,5.019101157591194,0.002686291252334517,3.670243749718892,0.2672231764451851,4.054779731679503,0.213841315950539),
    @(0.05183547313176007,0.008523908903537603,0.1042796413062703,5.021076329582144,0.002687955477581665,3.671754975806166,0.2673372526112132,4.021416905571357,0.2138913328255043),
    @(0.05179891215360755,0.008509059759044035,0.1042686760663276,5.021443851991179,0.002688234812156401,3.672035038292577,0.2673579110396105,4.015911190582131,0.2139008601940731),
    @(0.05205412879904969,0.008612120153498992,0.1043481590732469,5.019125143015984,0.00268631349627237,3.670262177193166,0.2672245998514384,4.054327499011436,0.2138419085754464),
    @(0.05321576171184716,0.009066555403791199,0.1047822012540678,5.014656274545047,0.002678268711259572,3.666683572950774,0.2668859601478104,4.228817017375093,0.2137593260269952),
    @(0.05561837089006616,0.009955337522342234,0.1059315538301924,5.02675532618818,0.002664039370454601,3.675034963709237,0.2671211744397795,4.588890451584518,0.2142350570981577),
    @(0.05745755474458747,0.01060603311495001,0.1069570207183759,5.048512641484791,0.002654516595829473,3.690653403633064,0.267844705351294,4.864446126579708,0.2149733500469821),
    @(0.05830992165117976,0.01090139217088648,0.1074627122233167,5.061238048405755,0.002650384305059333,3.699842826400342,0.2682940299432843,4.992219189850175,0.2153936873841715),
    @(0.05863491448386782,0.0110131303651535,0.1076598277186989,5.066466089072833,0.002648848041549813,3.703624349652756,0.2684815037769681,5.040953170658213,0.2155650071879833),
    @(0.05856482334630897,0.0109890705739204,0.1076171256171783,5.065321896370023,0.002649177636483333,3.702796482170712,0.2684403567475897,5.030441894406579,0.2155275701660884),
    @(0.05833661473774043,0.01091058718531457,0.1074788164674381,5.061659946612252,0.002650257344491242,3.700147877129538,0.2683091060783411,4.996221561880645,0.2154075385248362),
    @(0.0581971183283656,0.01086249936188821,0.1073948296904454,5.059470263705464,0.00265092240948696,3.698564876695087,0.2682309686595872,4.97530611103042,0.2153355975738478),
    @(0.05740216334542225,0.01058671641227704,0.1069247600351986,5.047738128556801,0.002654790654101085,3.690094964099501,0.2678177629230163,4.856144655810397,0.214947579892268),
    @(0.05691848022381407,0.01041735620274764,0.1066464163140424,5.041266944679307,0.002657214718053421,3.68543421250186,0.2675950864552661,4.783663750027699,0.2147311805669077),
    @(0.05664175975783792,0.01031988476855972,0.1064900111056239,5.037810910953112,0.002658627778983218,3.682949558753222,0.2674783188060061,4.742202503954843,0.2146146654703642),
    @(0.05654832259399001,0.01028687279471541,0.1064376893782395,5.036686369087562,0.002659109451860571,3.682141917197015,0.2674407245387371,4.728203546810505,0.2145765813022962),
    @(0.05696981610082474,0.01043539119956094,0.106675664619555,5.04192826277,0.002656954727352252,3.685910049224745,0.2676176199105385,4.791355877985723,0.2147533936502555),
    @(0.05840358519841971,0.01093364270294472,0.1075192887862251,5.062724423762063,0.002649939434860341,3.700917633336118,0.2683471870637817,5.006263421038682,0.2154424650634112),
    @(0.05935355835913469,0.01125864368636798,0.108103403126762,5.078702289077768,0.002645520842856693,3.71248534400398,0.2689250014406213,5.148753206760773,0.215963616227107),
    @(0.05884537020967429,0.01108524742902262,0.107788658117304,5.069955394760569,0.00264786396634566,3.706149808516713,0.2686073553568207,5.072517173682911,0.2156789897526252),
    @(0.05694660292611786,0.01042723790000366,0.1066624301892318,5.041628457706338,0.002657072208436799,3.685694316598159,0.2676073974908064,4.787877616037747,0.2147433265215781),
    @(0.05495526782093663,0.009715242460501372,0.1055887799553972,5.021233358128285,0.00266772438894185,3.671118696493082,0.2669610157205327,4.489555026962023,0.2140380927618963)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $col = $cols[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}
